# Updated symbol list on Thu Dec 15 22:52:31 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price ("D") column updates — these are stored as text in the workbook
# (numeric-looking strings), so force text format before assigning so
# Excel does not silently convert them to numbers.
$priceUpdates = @{
    "D2"  = "257.45"
    "D3"  = "22.79"
    "D4"  = "6.172"
    "D5"  = "0.06066"
    "D7"  = "3.455"
    "D8"  = "1.355"
    "D9"  = "0.7955"
    "D10" = "0.1579"
    "D11" = "0.08055"
    "D12" = "0.03348"
    "D13" = "0.03087"
    "D14" = "0.09291"
    "D15" = "3.926"
    "D16" = "0.001714"
    "D17" = "0.04842"
    "D18" = "0.0006147"
    "D19" = "0.006198"
    "D20" = "0.001100"
    "D21" = "0.003382"
    "D22" = "0.0001501"
    "D23" = "3.684"
    "D24" = "2.261"
    "D25" = "0.3357"
    "D27" = "0.0003019"
    "D40" = "0.04572"
    "D41" = "0.007146"
    "D42" = "0.1113"
    "D43" = "0.003132"
    "D44" = "0.009933"
    "D46" = "0.00005948"
    "D47" = "0.00000000751"
    "D48" = "0.7508"
    "D49" = "0.1062"
    "D50" = "0.00001502"
    "D51" = "0.01011"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Row 42 and 43 coin identities were swapped (CEJI <-> BKEXToken),
# with new price/volume-id data for each.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
